# fix: resolve #3 adding Restore tabs widnow popup handling mechanisim in diffrent stages
#
# The "Profile"/"Profile 2"/"Profile 4" column (A) previously had a
# duplicate/secondary column (B) holding "Profile 1"/"Profile 2"/"Profile 4"
# and "Default". Collapse the two-column layout back down to a single
# column: keep A2/A3 as "Profile 2"/"Profile 4" (the values that used to
# live in column B) and drop column B entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Profile 2"
$ws.Range("A3").Value = "Profile 4"

# Column B is no longer used - clear it out so the sheet's used range
# shrinks back down to a single column.
$ws.Range("B2:B3").ClearContents()

# Matches the updated selection recorded in the sheet view.
$ws.Range("D9").Select()
